$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new blank rows before the existing row 307, pushing the
# current rows 307-319 down to become rows 314-326. The inserted rows
# inherit the formatting (e.g. the date style on column D) of the row
# above, matching the target dimension change A1:T319 -> A1:T326.
$ws.Rows("307:313").Insert()

# Common values shared by every row in this data block.
$marketId   = 8
$market     = "Terminal La Palmera de La Serena"
$region     = "Coquimbo"
$codreg     = 4
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria  = "Durazno"
$unidad     = "`$/bins (400 kilos)"
$kgUnidad   = 400
$fecha      = 44585

# New row 307: Andross / Especial
$ws.Cells.Item(307, 1).Value = $marketId
$ws.Cells.Item(307, 2).Value = $market
$ws.Cells.Item(307, 3).Value = $region
$ws.Cells.Item(307, 4).Value = $fecha
$ws.Cells.Item(307, 5).Value = $codreg
$ws.Cells.Item(307, 6).Value = $tipo
$ws.Cells.Item(307, 7).Value = $productoId
$ws.Cells.Item(307, 8).Value = $producto
$ws.Cells.Item(307, 9).Value = $categoriaId
$ws.Cells.Item(307, 10).Value = $categoria
$ws.Cells.Item(307, 11).Value = "Andross"
$ws.Cells.Item(307, 12).Value = "Especial"
$ws.Cells.Item(307, 13).Value = 16
$ws.Cells.Item(307, 14).Value = 365000
$ws.Cells.Item(307, 15).Value = 370000
$ws.Cells.Item(307, 16).Value = 367500
$ws.Cells.Item(307, 17).Value = $unidad
$ws.Cells.Item(307, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(307, 19).Value = 919
$ws.Cells.Item(307, 20).Value = $kgUnidad

# New row 308: Andross / Primera
$ws.Cells.Item(308, 1).Value = $marketId
$ws.Cells.Item(308, 2).Value = $market
$ws.Cells.Item(308, 3).Value = $region
$ws.Cells.Item(308, 4).Value = $fecha
$ws.Cells.Item(308, 5).Value = $codreg
$ws.Cells.Item(308, 6).Value = $tipo
$ws.Cells.Item(308, 7).Value = $productoId
$ws.Cells.Item(308, 8).Value = $producto
$ws.Cells.Item(308, 9).Value = $categoriaId
$ws.Cells.Item(308, 10).Value = $categoria
$ws.Cells.Item(308, 11).Value = "Andross"
$ws.Cells.Item(308, 12).Value = "Primera"
$ws.Cells.Item(308, 13).Value = 20
$ws.Cells.Item(308, 14).Value = 325000
$ws.Cells.Item(308, 15).Value = 330000
$ws.Cells.Item(308, 16).Value = 327500
$ws.Cells.Item(308, 17).Value = $unidad
$ws.Cells.Item(308, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(308, 19).Value = 819
$ws.Cells.Item(308, 20).Value = $kgUnidad

# New row 309: Carson / Primera
$ws.Cells.Item(309, 1).Value = $marketId
$ws.Cells.Item(309, 2).Value = $market
$ws.Cells.Item(309, 3).Value = $region
$ws.Cells.Item(309, 4).Value = $fecha
$ws.Cells.Item(309, 5).Value = $codreg
$ws.Cells.Item(309, 6).Value = $tipo
$ws.Cells.Item(309, 7).Value = $productoId
$ws.Cells.Item(309, 8).Value = $producto
$ws.Cells.Item(309, 9).Value = $categoriaId
$ws.Cells.Item(309, 10).Value = $categoria
$ws.Cells.Item(309, 11).Value = "Carson"
$ws.Cells.Item(309, 12).Value = "Primera"
$ws.Cells.Item(309, 13).Value = 20
$ws.Cells.Item(309, 14).Value = 325000
$ws.Cells.Item(309, 15).Value = 330000
$ws.Cells.Item(309, 16).Value = 327500
$ws.Cells.Item(309, 17).Value = $unidad
$ws.Cells.Item(309, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(309, 19).Value = 819
$ws.Cells.Item(309, 20).Value = $kgUnidad

# New row 310: Carson / Segunda
$ws.Cells.Item(310, 1).Value = $marketId
$ws.Cells.Item(310, 2).Value = $market
$ws.Cells.Item(310, 3).Value = $region
$ws.Cells.Item(310, 4).Value = $fecha
$ws.Cells.Item(310, 5).Value = $codreg
$ws.Cells.Item(310, 6).Value = $tipo
$ws.Cells.Item(310, 7).Value = $productoId
$ws.Cells.Item(310, 8).Value = $producto
$ws.Cells.Item(310, 9).Value = $categoriaId
$ws.Cells.Item(310, 10).Value = $categoria
$ws.Cells.Item(310, 11).Value = "Carson"
$ws.Cells.Item(310, 12).Value = "Segunda"
$ws.Cells.Item(310, 13).Value = 20
$ws.Cells.Item(310, 14).Value = 265000
$ws.Cells.Item(310, 15).Value = 270000
$ws.Cells.Item(310, 16).Value = 267500
$ws.Cells.Item(310, 17).Value = $unidad
$ws.Cells.Item(310, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(310, 19).Value = 669
$ws.Cells.Item(310, 20).Value = $kgUnidad

# New row 311: Elegant Lady / Especial
$ws.Cells.Item(311, 1).Value = $marketId
$ws.Cells.Item(311, 2).Value = $market
$ws.Cells.Item(311, 3).Value = $region
$ws.Cells.Item(311, 4).Value = $fecha
$ws.Cells.Item(311, 5).Value = $codreg
$ws.Cells.Item(311, 6).Value = $tipo
$ws.Cells.Item(311, 7).Value = $productoId
$ws.Cells.Item(311, 8).Value = $producto
$ws.Cells.Item(311, 9).Value = $categoriaId
$ws.Cells.Item(311, 10).Value = $categoria
$ws.Cells.Item(311, 11).Value = "Elegant Lady"
$ws.Cells.Item(311, 12).Value = "Especial"
$ws.Cells.Item(311, 13).Value = 16
$ws.Cells.Item(311, 14).Value = 375000
$ws.Cells.Item(311, 15).Value = 380000
$ws.Cells.Item(311, 16).Value = 377500
$ws.Cells.Item(311, 17).Value = $unidad
$ws.Cells.Item(311, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(311, 19).Value = 944
$ws.Cells.Item(311, 20).Value = $kgUnidad

# New row 312: Elegant Lady / Primera
$ws.Cells.Item(312, 1).Value = $marketId
$ws.Cells.Item(312, 2).Value = $market
$ws.Cells.Item(312, 3).Value = $region
$ws.Cells.Item(312, 4).Value = $fecha
$ws.Cells.Item(312, 5).Value = $codreg
$ws.Cells.Item(312, 6).Value = $tipo
$ws.Cells.Item(312, 7).Value = $productoId
$ws.Cells.Item(312, 8).Value = $producto
$ws.Cells.Item(312, 9).Value = $categoriaId
$ws.Cells.Item(312, 10).Value = $categoria
$ws.Cells.Item(312, 11).Value = "Elegant Lady"
$ws.Cells.Item(312, 12).Value = "Primera"
$ws.Cells.Item(312, 13).Value = 16
$ws.Cells.Item(312, 14).Value = 335000
$ws.Cells.Item(312, 15).Value = 340000
$ws.Cells.Item(312, 16).Value = 337500
$ws.Cells.Item(312, 17).Value = $unidad
$ws.Cells.Item(312, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(312, 19).Value = 844
$ws.Cells.Item(312, 20).Value = $kgUnidad

# New row 313: Elegant Lady / Segunda
$ws.Cells.Item(313, 1).Value = $marketId
$ws.Cells.Item(313, 2).Value = $market
$ws.Cells.Item(313, 3).Value = $region
$ws.Cells.Item(313, 4).Value = $fecha
$ws.Cells.Item(313, 5).Value = $codreg
$ws.Cells.Item(313, 6).Value = $tipo
$ws.Cells.Item(313, 7).Value = $productoId
$ws.Cells.Item(313, 8).Value = $producto
$ws.Cells.Item(313, 9).Value = $categoriaId
$ws.Cells.Item(313, 10).Value = $categoria
$ws.Cells.Item(313, 11).Value = "Elegant Lady"
$ws.Cells.Item(313, 12).Value = "Segunda"
$ws.Cells.Item(313, 13).Value = 20
$ws.Cells.Item(313, 14).Value = 305000
$ws.Cells.Item(313, 15).Value = 310000
$ws.Cells.Item(313, 16).Value = 307500
$ws.Cells.Item(313, 17).Value = $unidad
$ws.Cells.Item(313, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 19).Value = 769
$ws.Cells.Item(313, 20).Value = $kgUnidad
